# Update countries & provincias Spain
# Applies the data refresh to the "Pais" sheet:
#  - Updates the "last updated" timestamp
#  - Refreshes case counts for several countries; the re-sort of the live
#    data feed causes a few adjacent rows to change which country they
#    display (their case totals crossed over), so those rows get the
#    country label shifted down one position along with the row's new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value2 = "Datos actualizados a 15 de Agosto de 2020 a las 09:21"

# --- Row 55/56: Kirguistan / Armenia swap places (Armenia overtakes) --
$ws.Range("A55").Value2 = "Armenia"
$ws.Range("B55").Value2 = 41495
$ws.Range("C55").Value2 = 196
$ws.Range("D55").Value2 = 34484
$ws.Range("E55").Value2 = 6194
$ws.Range("F55").Value2 = 0
$ws.Range("G55").Value2 = 3
$ws.Range("H55").Value2 = 817

$ws.Range("A56").Value2 = "Kirguistan"
$ws.Range("B56").Value2 = 41373
$ws.Range("C56").Value2 = 0
$ws.Range("D56").Value2 = 33592
$ws.Range("E56").Value2 = 6290
$ws.Range("F56").Value2 = 0
$ws.Range("G56").Value2 = 0
$ws.Range("H56").Value2 = 1491

# --- Row 73: El Salvador data refresh ---------------------------------
$ws.Range("D73").Value2 = 10474
$ws.Range("E73").Value2 = 11245

# --- Row 108: Hungria data refresh ------------------------------------
$ws.Range("B108").Value2 = 4877
$ws.Range("C108").Value2 = 24
$ws.Range("D108").Value2 = 3606
$ws.Range("E108").Value2 = 664

# --- Rows 145-147: Georgia overtakes Republica de Chipre & Letonia ----
$ws.Range("A145").Value2 = "Georgia"
$ws.Range("B145").Value2 = 1321
$ws.Range("C145").Value2 = 15
$ws.Range("D145").Value2 = 1088
$ws.Range("E145").Value2 = 216
$ws.Range("F145").Value2 = 0
$ws.Range("G145").Value2 = 0
$ws.Range("H145").Value2 = 17

$ws.Range("A146").Value2 = "Republica de Chipre"
$ws.Range("B146").Value2 = 1318
$ws.Range("C146").Value2 = 0
$ws.Range("D146").Value2 = 870
$ws.Range("E146").Value2 = 428
$ws.Range("F146").Value2 = 0
$ws.Range("G146").Value2 = 0
$ws.Range("H146").Value2 = 20

$ws.Range("A147").Value2 = "Letonia"
$ws.Range("B147").Value2 = 1308
$ws.Range("C147").Value2 = 0
$ws.Range("D147").Value2 = 1078
$ws.Range("E147").Value2 = 198
$ws.Range("F147").Value2 = 0
$ws.Range("G147").Value2 = 0
$ws.Range("H147").Value2 = 32

# --- Rows 213/214: Montserrat / Islas Malvinas swap places ------------
$ws.Range("A213").Value2 = "Islas Malvinas"
$ws.Range("D213").Value2 = 13
$ws.Range("H213").Value2 = 0

$ws.Range("A214").Value2 = "Montserrat"
$ws.Range("D214").Value2 = 12
$ws.Range("H214").Value2 = 1
